# Auto-generated COM-interop script applying the DailyStats update (po 15. 08. 2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H increments for rows 792-879 (PCR result lag correction, +1) ---
$hOnlyUpdates = @{
    792 = 465;
    793 = 416;
    794 = 420;
    795 = 424;
    796 = 408;
    797 = 370;
    798 = 345;
    799 = 332;
    800 = 291;
    801 = 298;
    802 = 298;
    803 = 288;
    804 = 267;
    805 = 262;
    806 = 252;
    807 = 237;
    808 = 235;
    809 = 243;
    810 = 242;
    811 = 233;
    812 = 222;
    813 = 209;
    814 = 196;
    815 = 192;
    816 = 195;
    817 = 196;
    818 = 178;
    819 = 173;
    820 = 173;
    821 = 157;
    822 = 153;
    823 = 160;
    824 = 164;
    825 = 151;
    826 = 156;
    827 = 152;
    828 = 142;
    829 = 149;
    830 = 155;
    831 = 145;
    832 = 143;
    833 = 162;
    834 = 166;
    835 = 156;
    836 = 165;
    837 = 182;
    838 = 186;
    839 = 214;
    840 = 219;
    841 = 236;
    842 = 244;
    843 = 259;
    844 = 276;
    845 = 305;
    846 = 311;
    847 = 317;
    848 = 329;
    849 = 314;
    850 = 325;
    851 = 347;
    852 = 333;
    853 = 343;
    854 = 331;
    855 = 331;
    856 = 320;
    857 = 336;
    858 = 352;
    859 = 379;
    860 = 362;
    861 = 377;
    862 = 395;
    863 = 386;
    864 = 424;
    865 = 474;
    866 = 500;
    867 = 533;
    868 = 545;
    869 = 563;
    870 = 590;
    871 = 641;
    872 = 697;
    873 = 719;
    874 = 739;
    875 = 757;
    876 = 734;
    877 = 710;
    878 = 718;
    879 = 745
}
foreach ($row in $hOnlyUpdates.Keys) {
    $ws.Cells.Item($row, 8).Value = $hOnlyUpdates[$row]
}

# --- Multi-column updates for rows 880, 881-890 (revised F/G/H figures) ---
$multiColUpdates = @{
    880 = @{'F'=4935; 'G'=572; 'H'=718}
    881 = @{'H'=696}
    882 = @{'H'=648}
    883 = @{'H'=652}
    884 = @{'F'=3208; 'G'=324; 'H'=580}
    885 = @{'F'=1080; 'H'=604}
    886 = @{'F'=1450; 'H'=640}
    887 = @{'F'=4872; 'G'=496; 'H'=608}
    888 = @{'F'=3201; 'G'=285; 'H'=582}
    889 = @{'F'=3049; 'G'=273; 'H'=580}
    890 = @{'F'=3349; 'G'=234; 'H'=581}
}
foreach ($row in $multiColUpdates.Keys) {
    $cols = $multiColUpdates[$row]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$row").Value = $cols[$colLetter]
    }
}

# --- New rows 891-893 (daily stats through 2022-08-14) ---
$newRows = @{
    891 = @{'A'=44785; 'B'=1827404; 'C'=1912; 'D'=776; 'E'=20296; 'F'=2379; 'G'=283; 'H'=580}
    892 = @{'A'=44786; 'B'=1827711; 'C'=750; 'D'=307; 'E'=20303; 'F'=683; 'G'=44; 'H'=580}
    893 = @{'A'=44787; 'B'=1827804; 'C'=316; 'D'=93; 'E'=20306; 'F'=463; 'G'=54; 'H'=580}
}
foreach ($row in ($newRows.Keys | Sort-Object)) {
    $cols = $newRows[$row]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$row").Value = $cols[$colLetter]
    }
}

Write-Output "Applied DailyStats update through 2022-08-14 (rows 792-893)."
